$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = 5.790099999999998
$ws.Range("E4").Value  = 13.7638
$ws.Range("E5").Value  = 13.233
$ws.Range("B6").Value  = 9.277300000000006
$ws.Range("B7").Value  = 6.638799999999995
$ws.Range("E8").Value  = 14.09509999999999
$ws.Range("B16").Value = 8.935600000000004
$ws.Range("E16").Value = 12.69170000000001
$ws.Range("B20").Value = 5.803199999999997
$ws.Range("E22").Value = 12.3654
